{"js": "// Replace the date and the multiplication problems throughout the document.\n// All source strings in this worksheet are unique, so a straightforward\n// search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-03-17 Sunday\", \"2024-03-18 Monday\"],\n  [\"771\u00d74=\", \"228\u00d76=\"],\n  [\"592\u00d76=\", \"924\u00d76=\"],\n  [\"525\u00d74=\", \"732\u00d74=\"],\n  [\"865\u00d72=\", \"146\u00d77=\"],\n  [\"835\u00d77=\", \"370\u00d72=\"],\n  [\"306\u00d72=\", \"374\u00d73=\"],\n  [\"372\u00d79=\", \"476\u00d75=\"],\n  [\"976\u00d74=\", \"623\u00d73=\"],\n  [\"736\u00d72=\", \"396\u00d73=\"],\n  [\"865\u00d74=\", \"510\u00d72=\"],\n  [\"221\u00d73=\", \"221\u00d75=\"],\n  [\"878\u00d79=\", \"432\u00d72=\"],\n  [\"354\u00d78=\", \"763\u00d77=\"],\n  [\"500\u00d75=\", \"559\u00d78=\"],\n  [\"294\u00d74=\", \"267\u00d75=\"],\n  [\"395\u00d76=\", \"285\u00d79=\"],\n  [\"127\u00d73=\", \"226\u00d76=\"],\n  [\"300\u00d75=\", \"908\u00d78=\"],\n  [\"199\u00d75=\", \"710\u00d75=\"],\n  [\"407\u00d76=\", \"474\u00d75=\"],\n  [\"233\u00d75=\", \"883\u00d73=\"],\n  [\"382\u00d74=\", \"175\u00d74=\"],\n  [\"685\u00d73=\", \"740\u00d72=\"],\n  [\"176\u00d75=\", \"303\u00d72=\"],\n  [\"404\u00d77=\", \"376\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the multiplication problems throughout the document.\n# All source strings in this worksheet are unique, so a straightforward\n# Find/Replace pass per pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-17 Sunday\", \"2024-03-18 Monday\"),\n    @(\"771\u00d74=\", \"228\u00d76=\"),\n    @(\"592\u00d76=\", \"924\u00d76=\"),\n    @(\"525\u00d74=\", \"732\u00d74=\"),\n    @(\"865\u00d72=\", \"146\u00d77=\"),\n    @(\"835\u00d77=\", \"370\u00d72=\"),\n    @(\"306\u00d72=\", \"374\u00d73=\"),\n    @(\"372\u00d79=\", \"476\u00d75=\"),\n    @(\"976\u00d74=\", \"623\u00d73=\"),\n    @(\"736\u00d72=\", \"396\u00d73=\"),\n    @(\"865\u00d74=\", \"510\u00d72=\"),\n    @(\"221\u00d73=\", \"221\u00d75=\"),\n    @(\"878\u00d79=\", \"432\u00d72=\"),\n    @(\"354\u00d78=\", \"763\u00d77=\"),\n    @(\"500\u00d75=\", \"559\u00d78=\"),\n    @(\"294\u00d74=\", \"267\u00d75=\"),\n    @(\"395\u00d76=\", \"285\u00d79=\"),\n    @(\"127\u00d73=\", \"226\u00d76=\"),\n    @(\"300\u00d75=\", \"908\u00d78=\"),\n    @(\"199\u00d75=\", \"710\u00d75=\"),\n    @(\"407\u00d76=\", \"474\u00d75=\"),\n    @(\"233\u00d75=\", \"883\u00d73=\"),\n    @(\"382\u00d74=\", \"175\u00d74=\"),\n    @(\"685\u00d73=\", \"740\u00d72=\"),\n    @(\"176\u00d75=\", \"303\u00d72=\"),\n    @(\"404\u00d77=\", \"376\u00d79=\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
